$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "306.02"
$ws.Range("E2").Value = "0.90%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "35.98"
$ws.Range("E3").Value = "0.78%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.012"
$ws.Range("E4").Value = "-1.19%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08093"
$ws.Range("E5").Value = "0.31%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "1.921"
$ws.Range("E6").Value = "0.07%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "4.142"
$ws.Range("E7").Value = "2.23%"

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "7.875"
$ws.Range("E8").Value = "1.11%"

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9310"
$ws.Range("E9").Value = "0.19%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1254"
$ws.Range("E10").Value = "-15.98%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1902"
$ws.Range("E11").Value = "0.22%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09240"
$ws.Range("E12").Value = "3.02%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.72%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09915"
$ws.Range("E14").Value = "0.73%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001419"
$ws.Range("E15").Value = "0.03%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006405"
$ws.Range("E16").Value = "11.29%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.99%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.285"
$ws.Range("E18").Value = "10.15%"

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3437"
$ws.Range("E19").Value = "-0.13%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "5.201"
$ws.Range("E20").Value = "3.42%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1284"
$ws.Range("E21").Value = "-1.15%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2528"
$ws.Range("E22").Value = "1.42%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04415"
$ws.Range("E23").Value = "-1.97%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").Value = "2.08%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004729"
$ws.Range("E25").Value = "-1.67%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").Value = "6.12%"

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003125"
$ws.Range("E27").Value = "3.68%"

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01969"
$ws.Range("E39").Value = "5.12%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05243"
$ws.Range("E40").Value = "9.51%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007543"
$ws.Range("E41").Value = "2.97%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01017"
$ws.Range("E42").Value = "-4.02%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1375"
$ws.Range("E43").Value = "2.34%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002103"
$ws.Range("E44").Value = "-0.06%"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "9.82%"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006379"
$ws.Range("E46").Value = "2.69%"

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "0.16%"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "65.22"
$ws.Range("E48").Value = "0.86%"

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001657"
$ws.Range("E49").Value = "-0.13%"

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").Value = "0.16%"

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").Value = "0.16%"
